$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.085.05"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "2.105.80"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.17%  "
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5156"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("D13").Value = "2.113.43"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.756"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06667"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.283"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "30.178.61"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.353"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "2.360.37"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.562"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.184"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.644"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.270"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.983"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.921"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02585"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06853"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2316"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6844"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.254"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6427"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("E47").Value = "  +3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.662"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "83.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07245"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.88%  "
